$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.749289631843567
$ws.Range("B1").Value = 2.033302545547485
$ws.Range("C1").Value = 2.578598499298096
$ws.Range("D1").Value = 4.2105393409729
$ws.Range("E1").Value = 2.939725399017334
